$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95; existing rows 95-153 shift down to 96-154.
$ws.Rows(95).Insert()

# Populate the newly inserted row 95 with its data.
$ws.Cells.Item(95, 1).Value = 5
$ws.Cells.Item(95, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(95, 3).Value = "Maule"
$ws.Cells.Item(95, 4).Value = 45001
$ws.Cells.Item(95, 5).Value = 7
$ws.Cells.Item(95, 6).Value = 100112001
$ws.Cells.Item(95, 7).Value = "Berenjena"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 250
$ws.Cells.Item(95, 11).Value = 8000
$ws.Cells.Item(95, 12).Value = 8000
$ws.Cells.Item(95, 13).Value = 8000
$ws.Cells.Item(95, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(95, 15).Value = "Región del Maule"
$ws.Cells.Item(95, 16).Value = 160
$ws.Cells.Item(95, 17).Value = 50
$ws.Cells.Item(95, 18).Value = "Hortaliza"
